$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-24 Wednesday" "2025-12-25 Thursday"

Replace-Text "907÷5=" "397÷4="
Replace-Text "118÷5=" "911÷2="
Replace-Text "725÷4=" "150÷9="
Replace-Text "532÷7=" "826÷2="
Replace-Text "101÷8=" "807÷9="
Replace-Text "766÷3=" "491÷5="
Replace-Text "821÷9=" "744÷9="
Replace-Text "179÷7=" "130÷2="
Replace-Text "388÷9=" "782÷2="
Replace-Text "696÷6=" "539÷9="
Replace-Text "822÷2=" "552÷2="
Replace-Text "759÷8=" "213÷5="
Replace-Text "766÷9=" "270÷2="
Replace-Text "906÷6=" "400÷3="
Replace-Text "668÷8=" "913÷6="
Replace-Text "962÷5=" "415÷2="
Replace-Text "263÷3=" "821÷9="
Replace-Text "834÷6=" "963÷4="
Replace-Text "506÷2=" "426÷9="
Replace-Text "995÷8=" "763÷5="
Replace-Text "428÷4=" "246÷5="
Replace-Text "711÷5=" "661÷9="
Replace-Text "499÷3=" "796÷4="
Replace-Text "966÷5=" "991÷2="
Replace-Text "792÷9=" "928÷3="
